# Interim commit - changes to autonomous variables speed

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1 data edits ---

# Row 11: new cells I11, J11
$sheet1.Range("I11").Value = 19.5
$sheet1.Range("J11").Value = "skewed to right"

# Row 12: new cells A12, B12, C12, D12 (D12 is a formula)
$sheet1.Range("A12").Value = 11
$sheet1.Range("B12").Value = 0.5
$sheet1.Range("C12").Value = 0.000069999999999999994
$sheet1.Range("D12").Formula = "=20/4"

# --- View / selection state ---

$sheet1.Activate()
$excel.ActiveWindow.Zoom = 125
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$sheet1.Range("C12").Select()

$wb.Save()
